$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The special-character facet-encoding test (row 5) was removed from the
# QA test table ("Move to NIEM Releases-specific tests... doesn't apply
# to user input"). Deleting the whole row shifts every row below it up
# by one, which is exactly what the table/worksheet diff shows.
$ws.Rows.Item(5).Delete()

# Restore the active cell to A2 (matches the saved selection in the
# edited workbook).
$ws.Range("A2").Select()
